$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(43966, 56, 16, 155, 40),
    @(43967, 54, 15, 151, 38),
    @(43968, 51, 10, 138, 34),
    @(43969, 55, 4, 141, 39),
    @(43970, 54, 16, 132, 38),
    @(43971, 49, 11, 126, 35),
    @(43972, 52, 12, 109, 36)
)

$row = 34
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

# Match the author's final view state: scrolled down a bit further and the
# active selection moved on to the next blank row below the new data.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C41").Select()
